$wb = $excel.ActiveWorkbook

# Update the FD sheet inputs
$wsFD = $wb.Worksheets.Item("FD")
$wsFD.Range("B2").Value = 1000
$wsFD.Range("B4").Value = 0
$wsFD.Range("B6").Value = 1000

# Update selection on FD sheet to C11
$wsFD.Range("C11").Select()

# Make FD the active sheet / active tab
$wsFD.Activate()

$excel.Calculate()
